$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($ws, $addr, $text) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
}

Set-CellText $ws "D2" "302.78"
Set-CellText $ws "E2" "1.89%"
Set-CellText $ws "E3" "0.80%"
Set-CellText $ws "D4" "5.162"
Set-CellText $ws "E4" "1.13%"
Set-CellText $ws "D5" "0.07913"
Set-CellText $ws "E5" "5.68%"
Set-CellText $ws "D6" "2.338"
Set-CellText $ws "E6" "36.70%"
Set-CellText $ws "D7" "7.948"
Set-CellText $ws "E7" "2.75%"
Set-CellText $ws "E8" "1.87%"
Set-CellText $ws "D9" "0.9066"
Set-CellText $ws "E9" "-2.73%"
Set-CellText $ws "D10" "0.1736"
Set-CellText $ws "E10" "2.70%"
Set-CellText $ws "D11" "0.07385"
Set-CellText $ws "E11" "3.42%"
Set-CellText $ws "D12" "0.08163"
Set-CellText $ws "E12" "2.35%"
Set-CellText $ws "D13" "0.03104"
Set-CellText $ws "E13" "2.62%"
Set-CellText $ws "D14" "0.09940"
Set-CellText $ws "E14" "0.40%"
Set-CellText $ws "D15" "0.001517"
Set-CellText $ws "E15" "1.92%"
Set-CellText $ws "D16" "0.006044"
Set-CellText $ws "E16" "-3.92%"
Set-CellText $ws "D17" "3.499"
Set-CellText $ws "E17" "1.31%"
Set-CellText $ws "E18" "0.66%"
Set-CellText $ws "E19" "-1.27%"
Set-CellText $ws "E20" "0.92%"
Set-CellText $ws "D21" "4.684"
Set-CellText $ws "E21" "2.52%"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-CellText $ws "D22" "0.1648"
Set-CellText $ws "E22" "5.74%"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-CellText $ws "D23" "0.04663"
Set-CellText $ws "D24" "0.001265"
Set-CellText $ws "E24" "3.83%"
Set-CellText $ws "D25" "0.004514"
Set-CellText $ws "E25" "1.87%"
Set-CellText $ws "E26" "3.73%"
Set-CellText $ws "D27" "0.0002741"
Set-CellText $ws "E27" "46.06%"
Set-CellText $ws "D39" "0.01939"
Set-CellText $ws "E39" "14.54%"
Set-CellText $ws "D40" "0.04574"
Set-CellText $ws "E40" "2.46%"
Set-CellText $ws "D41" "0.007303"
Set-CellText $ws "E41" "3.36%"
Set-CellText $ws "D42" "0.1360"
Set-CellText $ws "E42" "2.36%"
Set-CellText $ws "D43" "0.002249"
Set-CellText $ws "E43" "9.12%"
Set-CellText $ws "E44" "-5.14%"
Set-CellText $ws "D45" "0.00006472"
Set-CellText $ws "E45" "7.98%"
Set-CellText $ws "E46" "-0.06%"
Set-CellText $ws "E48" "15.31%"
Set-CellText $ws "D49" "0.00002099"
Set-CellText $ws "E49" "-0.06%"
Set-CellText $ws "D50" "0.0001999"
Set-CellText $ws "E50" "0.01%"
